# feat: add 2022-Q1 data
#
# Before: sheets = 2021-Q2, 2021-Q4, 总计
# After:  sheets = 2021-Q2, 2021-Q4, 2022-Q1, 总计
#
# The old "总计" (grand-total) sheet is renamed to "2022-Q1" and its
# content is replaced with the new quarter's per-fund holdings row
# (same shape as the 2021-Q4 sheet). A brand-new "总计" sheet is
# inserted right after it, holding the updated roll-up table (with the
# new 2022-Q1 row prepended).

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# ------------------------------------------------------------------
# 1) Build the brand-new "总计" sheet FIRST (while the old totals
#    sheet still holds its original header/row formatting to copy
#    from), positioned immediately after the current "总计" sheet.
# ------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add($null, $oldTotal)
$newTotal.Name = "TEMP_NEW_TOTAL"

# Match page margins used by the rest of the workbook's sheets.
$newTotal.PageSetup.LeftMargin = 54
$newTotal.PageSetup.RightMargin = 54
$newTotal.PageSetup.TopMargin = 72
$newTotal.PageSetup.BottomMargin = 72
$newTotal.PageSetup.HeaderMargin = 36
$newTotal.PageSetup.FooterMargin = 36

# Header row (B1:D1) — copy text + style from the old totals sheet.
$oldTotal.Range("B1:D1").Copy($newTotal.Range("B1"))

# Row 2: new 2022-Q1 roll-up entry (prepended).
$oldTotal.Range("A2").Copy($newTotal.Range("A2"))
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 0

# Row 3: previous 2021-Q4 roll-up entry (shifted down one row).
$oldTotal.Range("A2").Copy($newTotal.Range("A3"))
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 1
$newTotal.Range("D3").Value = 0.01

# Row 4: previous 2021-Q2 roll-up entry (shifted down one row).
$oldTotal.Range("A2").Copy($newTotal.Range("A4"))
$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q2"
$newTotal.Range("C4").Value = 1
$newTotal.Range("D4").Value = 0.01

# ------------------------------------------------------------------
# 2) Turn the old "总计" sheet into the new "2022-Q1" quarter sheet,
#    mirroring the 2021-Q4 sheet's layout (same headers, same fund),
#    with this quarter's figures.
# ------------------------------------------------------------------
$oldTotal.Cells.Clear()
$oldTotal.Name = "2022-Q1"

# Header row + fund code/name, copied straight from 2021-Q4.
$q4.Range("B1:H1").Copy($oldTotal.Range("B1"))
$q4.Range("A2:C2").Copy($oldTotal.Range("A2"))

# Quarter-specific figures (kept as text, matching the source data's
# original typing, same as every other quarter sheet).
$oldTotal.Range("D2").NumberFormat = "@"
$oldTotal.Range("D2").Value = "0.33"
$oldTotal.Range("D2").Style = "Normal"

$oldTotal.Range("E2").NumberFormat = "@"
$oldTotal.Range("E2").Value = "37.77"
$oldTotal.Range("E2").Style = "Normal"

$oldTotal.Range("F2").NumberFormat = "@"
$oldTotal.Range("F2").Value = "0.86"
$oldTotal.Range("F2").Style = "Normal"

$oldTotal.Range("G2").NumberFormat = "@"
$oldTotal.Range("G2").Value = "0.0028"
$oldTotal.Range("G2").Style = "Normal"

$oldTotal.Range("H2").Value = 10

# ------------------------------------------------------------------
# 3) Finally, rename the new roll-up sheet to "总计" (it now lands
#    right after "2022-Q1" in tab order, as intended).
# ------------------------------------------------------------------
$newTotal.Name = "总计"

# Restore the originally active/selected sheet (adding sheets shifts
# focus onto them in Excel, but the workbook's own selection state is
# otherwise unrelated to this data edit).
$q2 = $wb.Worksheets.Item("2021-Q2")
$q2.Activate()
